# [Update] CDUN Agregar Hostel Worker
# Updates the "Curso normal (Básico)" narrative steps table: rewords three
# steps to explicitly name the Admin actor, and re-numbers/re-flows the
# six steps from a 2-column x 3-row layout into a layout where row 18's
# right-hand cell is emptied and a new row 20 is used for step 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1 (B17/C17) -------------------------------------------------
$ws.Range("C17").Value = "Actor Admin: Seleciona agregar Hostel Worker"

# --- Step 2: move "El sistema muestra un formulario..." from C18 to G17,
#     renumber F17 3 -> 2 ------------------------------------------------
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = "El sistema muestra un formulario para ingresar los datos."

# --- Step 3 (B18/C18): renumber 2 -> 3, reword -------------------------
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = "Actor Admin: Ingresa los datos personales(nombre, rut, dirección, telefono, fecha de nacimiento) del  Hostel Worker."

# --- Former step 4 (F18/G18) removed: clear the now-unused right cell --
$ws.Range("F18").ClearContents()
$ws.Range("G18:I18").ClearContents()

# --- Step 4 (B19/C19): renumber 5 -> 4, reword --------------------------
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = "Actor Admin: ingresa los un nombre de usuario y contraseña para el HostelWorker"

# --- Step 5: move "El sistema valida los datos..." from C19 to G19,
#     renumber F19 6 -> 5 -------------------------------------------------
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = "El sistema valida los datos y los ingresa al sistema."

# --- Step 6: now lives in a new row 20 (was G19) ------------------------
$ws.Range("B20").Value = 6
$ws.Range("C20").Value = "Fin del caso de uso"

# --- Apply matching styles to the new row 20 cells, mirroring row 19 ---
$ws.Range("B20").Style = $ws.Range("B19").Style
$ws.Range("C20").Style = $ws.Range("C19").Style

# --- Restore the bottom-border look of row 19/20 (row 20 is now the
#     closing, thick-bottom-bordered row of the table) ------------------
$ws.Range("B20:I20").Borders.Item(9).LineStyle = 1
$ws.Range("B20:I20").Borders.Item(9).Weight = -4138
$ws.Range("B19:I19").Borders.Item(9).LineStyle = -4142

# --- Final cursor / view position, matching the author's saved state ---
$ws.Range("J26").Select()
$excel.ActiveWindow.ScrollRow = 4
